$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append a new row of data (row 9) below the existing table (A1:H8).
# Columns C and F hold date-like text ("2025-01-08", "2025-01-15") that
# Excel would otherwise auto-convert into date serial numbers, so those
# two cells are forced to text first and then restored to the default
# "Normal" style (no explicit number format) to match plain string cells
# elsewhere in the sheet.

$ws.Range("A9").Value = "DOC-1735824110933"
$ws.Range("B9").Value = "sjfdhfs"

$ws.Range("C9").NumberFormat = "@"
$ws.Range("C9").Value = "2025-01-08"
$ws.Range("C9").Style = "Normal"

$ws.Range("D9").Value = "12:21"
$ws.Range("E9").Value = "pdf"

$ws.Range("F9").NumberFormat = "@"
$ws.Range("F9").Value = "2025-01-15"
$ws.Range("F9").Style = "Normal"

$ws.Range("G9").Value = "14:32"
$ws.Range("H9").Value = "aefa"
